$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.08244627342756239
$ws.Range("D2").Value = 0.04072730343426301
$ws.Range("E2").Value = 0.1147618466391851
$ws.Range("F2").Value = 4.212240608295616
$ws.Range("G2").Value = 3.614556738805618
$ws.Range("H2").Value = 2.592033047395887
$ws.Range("I2").Value = 3.070072893264268
$ws.Range("J2").Value = 0.243591450515332
$ws.Range("L2").Value = 0.2099675411600899
$ws.Range("N2").Value = 2.873581687696401
$ws.Range("C3").Value = 0.08274999149471007
$ws.Range("D3").Value = 0.03994293178534036
$ws.Range("E3").Value = 0.1151013490423587
$ws.Range("F3").Value = 4.132352122187086
$ws.Range("G3").Value = 3.517155678699027
$ws.Range("H3").Value = 2.555584062300682
$ws.Range("I3").Value = 3.009488150405645
$ws.Range("J3").Value = 0.242746324526145
$ws.Range("L3").Value = 0.2101991750566796
$ws.Range("N3").Value = 2.562605684679681
$ws.Range("C4").Value = 0.08295171857653827
$ws.Range("D4").Value = 0.03945342185560818
$ws.Range("E4").Value = 0.115350064938168
$ws.Range("F4").Value = 4.085917180444056
$ws.Range("G4").Value = 3.459749556915909
$ws.Range("H4").Value = 2.534735381654457
$ws.Range("I4").Value = 2.974211608595198
$ws.Range("J4").Value = 0.2423468978063212
$ws.Range("L4").Value = 0.2104273127012846
$ws.Range("N4").Value = 2.371325805375818
$ws.Range("C5").Value = 0.08303776512076411
$ws.Range("D5").Value = 0.03925191464412237
$ws.Range("E5").Value = 0.1154615455231749
$ws.Range("F5").Value = 4.067648321696311
$ws.Range("G5").Value = 3.436954209929382
$ws.Range("H5").Value = 2.526622299106918
$ws.Range("I5").Value = 2.960316374123423
$ws.Range("J5").Value = 0.2422140927897019
$ws.Range("L5").Value = 0.2105418663036218
$ws.Range("N5").Value = 2.293303068607429
$ws.Range("C6").Value = 0.08305228532074871
$ws.Range("D6").Value = 0.03921833085809467
$ws.Range("E6").Value = 0.1154806685285301
$ws.Range("F6").Value = 4.064654153319253
$ws.Range("G6").Value = 3.433205039765141
$ws.Range("H6").Value = 2.525298195254862
$ws.Range("I6").Value = 2.958037992923252
$ws.Range("J6").Value = 0.242193848104499
$ws.Range("L6").Value = 0.2105621910660531
$ws.Range("N6").Value = 2.280343261403573
$ws.Range("C7").Value = 0.082952863465799
$ws.Range("D7").Value = 0.0394507125239798
$ws.Range("E7").Value = 0.1153515273962142
$ws.Range("F7").Value = 4.085668158771142
$ws.Range("G7").Value = 3.459439715698295
$ws.Range("H7").Value = 2.534624418063373
$ws.Range("I7").Value = 2.974022272345664
$ws.Range("J7").Value = 0.2423449855294635
$ws.Range("L7").Value = 0.2104287702335412
$ws.Range("N7").Value = 2.370273851392596
$ws.Range("C8").Value = 0.08254783659764797
$ws.Range("D8").Value = 0.04045846743854398
$ws.Range("E8").Value = 0.1148705526880338
$ws.Range("F8").Value = 4.184149277594287
$ws.Range("G8").Value = 3.580471956284896
$ws.Range("H8").Value = 2.579146356522813
$ws.Range("I8").Value = 3.048782211994336
$ws.Range("J8").Value = 0.2432752020111408
$ws.Range("L8").Value = 0.2100295655396565
$ws.Range("N8").Value = 2.766433886209825
$ws.Range("C9").Value = 0.08187417311815537
$ws.Range("D9").Value = 0.04237388412274967
$ws.Range("E9").Value = 0.1142467577314488
$ws.Range("F9").Value = 4.398246425471086
$ws.Range("G9").Value = 3.837093253837224
$ws.Range("H9").Value = 2.678706758664759
$ws.Range("I9").Value = 3.21080252247117
$ws.Range("J9").Value = 0.2460515404226769
$ws.Range("L9").Value = 0.2099294761749846
$ws.Range("N9").Value = 3.540180268007646
$ws.Range("C10").Value = 0.08145228262232074
$ws.Range("D10").Value = 0.04374684716379562
$ws.Range("E10").Value = 0.1139832305969293
$ws.Range("F10").Value = 4.568651750589595
$ws.Range("G10").Value = 4.03775519853076
$ws.Range("H10").Value = 2.759481704697009
$ws.Range("I10").Value = 3.339479224269553
$ws.Range("J10").Value = 0.2486782488693677
$ws.Range("L10").Value = 0.2102740097250475
$ws.Range("N10").Value = 4.10623028343673
$ws.Range("C11").Value = 0.08127612052366118
$ws.Range("D11").Value = 0.0443646086006666
$ws.Range("E11").Value = 0.1139056729405752
$ws.Range("F11").Value = 4.649092480800192
$ws.Range("G11").Value = 4.131757295952866
$ws.Range("H11").Value = 2.797920588990053
$ws.Range("I11").Value = 3.400165733711361
$ws.Range("J11").Value = 0.2500021465126068
$ws.Range("L11").Value = 0.2105219793652182
$ws.Range("N11").Value = 4.363110593465422
$ws.Range("C12").Value = 0.08121167117837658
$ws.Range("D12").Value = 0.0445976086629436
$ws.Range("E12").Value = 0.1138823922789758
$ws.Range("F12").Value = 4.679979194736944
$ws.Range("G12").Value = 4.16775107815937
$ws.Range("H12").Value = 2.81272278311252
$ws.Range("I12").Value = 3.423459690805061
$ws.Range("J12").Value = 0.2505221380717373
$ws.Range("L12").Value = 0.2106290343934205
$ws.Range("N12").Value = 4.460285735714251
$ws.Range("C13").Value = 0.08122545111211998
$ws.Range("D13").Value = 0.04454746867643777
$ws.Range("E13").Value = 0.1138871353343482
$ws.Range("F13").Value = 4.673308175592808
$ws.Range("G13").Value = 4.159981407583018
$ws.Range("H13").Value = 2.80952387217917
$ws.Range("I13").Value = 3.418428921644562
$ws.Range("J13").Value = 0.2504093168564765
$ws.Range("L13").Value = 0.2106053925704003
$ws.Range("N13").Value = 4.439361943450422
$ws.Range("C14").Value = 0.08127077299166885
$ws.Range("D14").Value = 0.04438379608431475
$ws.Range("E14").Value = 0.1139036355800211
$ws.Range("F14").Value = 4.651624989708125
$ws.Range("G14").Value = 4.13471052956362
$ws.Range("H14").Value = 2.799133424391869
$ws.Range("I14").Value = 3.402075839405313
$ws.Range("J14").Value = 0.2500445519021284
$ws.Range("L14").Value = 0.2105305230116628
$ws.Range("N14").Value = 4.371107314139522
$ws.Range("C15").Value = 0.08129882800331067
$ws.Range("D15").Value = 0.04428342179733136
$ws.Range("E15").Value = 0.1139145354881883
$ws.Range("F15").Value = 4.638398991914016
$ws.Range("G15").Value = 4.119283301500673
$ws.Range("H15").Value = 2.792801120007141
$ws.Range("I15").Value = 3.392100025023012
$ws.Range("J15").Value = 0.2498235563125135
$ws.Range("L15").Value = 0.2104863773736199
$ws.Range("N15").Value = 4.329286057409945
$ws.Range("C16").Value = 0.08146411175550128
$ws.Range("D16").Value = 0.04370634195018042
$ws.Range("E16").Value = 0.1139891510304984
$ws.Range("F16").Value = 4.563454017970543
$ws.Range("G16").Value = 4.031667191158533
$ws.Range("H16").Value = 2.757003943429595
$ws.Range("I16").Value = 3.335556838464726
$ws.Range("J16").Value = 0.2485943341446415
$ws.Range("L16").Value = 0.210259643173984
$ws.Range("N16").Value = 4.089429168003562
$ws.Range("C17").Value = 0.08156953925634802
$ws.Range("D17").Value = 0.04335061133653895
$ws.Range("E17").Value = 0.1140457667777479
$ws.Range("F17").Value = 4.518229954706555
$ws.Range("G17").Value = 3.978618488073437
$ws.Range("H17").Value = 2.735479248933302
$ws.Range("I17").Value = 3.301423122165744
$ws.Range("J17").Value = 0.2478733625170548
$ws.Range("L17").Value = 0.2101439414016752
$ws.Range("N17").Value = 3.94211849063862
$ws.Range("C18").Value = 0.08163166196579752
$ws.Range("D18").Value = 0.04314536105426114
$ws.Range("E18").Value = 0.1140823140050671
$ws.Range("F18").Value = 4.492493022993926
$ws.Range("G18").Value = 3.94836192007358
$ws.Range("H18").Value = 2.723258095246081
$ws.Range("I18").Value = 3.281992512070104
$ws.Range("J18").Value = 0.2474708134321943
$ws.Range("L18").Value = 0.2100859785397233
$ws.Range("N18").Value = 3.857331695637754
$ws.Range("C19").Value = 0.08165295066008227
$ws.Range("D19").Value = 0.04307575493767501
$ws.Range("E19").Value = 0.1140953723362355
$ws.Range("F19").Value = 4.483825969884464
$ws.Range("G19").Value = 3.938161289513289
$ws.Range("H19").Value = 2.719147500196982
$ws.Range("I19").Value = 3.275448249123656
$ws.Range("J19").Value = 0.2473365980628515
$ws.Range("L19").Value = 0.2100678269111356
$ws.Range("N19").Value = 3.828614786364199
$ws.Range("C20").Value = 0.08155816282083705
$ws.Range("D20").Value = 0.04338854573108364
$ws.Range("E20").Value = 0.1140393276551848
$ws.Range("F20").Value = 4.523015659545507
$ws.Range("G20").Value = 3.984239111306977
$ws.Range("H20").Value = 2.737754085559743
$ws.Range("I20").Value = 3.305035758849954
$ws.Range("J20").Value = 0.2479488544886337
$ws.Range("L20").Value = 0.2101553692209777
$ws.Range("N20").Value = 3.95780600327754
$ws.Range("C21").Value = 0.08125739958836498
$ws.Range("D21").Value = 0.04443189563777494
$ws.Range("E21").Value = 0.1138986237848236
$ws.Range("F21").Value = 4.657982274591575
$ws.Range("G21").Value = 4.142122368998116
$ws.Range("H21").Value = 2.802178646993866
$ws.Range("I21").Value = 3.406870600930148
$ws.Range("J21").Value = 0.2501511848730971
$ws.Range("L21").Value = 0.2105521567295483
$ws.Range("N21").Value = 4.391158149571083
$ws.Range("C22").Value = 0.08107400013787647
$ws.Range("D22").Value = 0.04510837003673629
$ws.Range("E22").Value = 0.1138421573596489
$ws.Range("F22").Value = 4.748673783272295
$ws.Range("G22").Value = 4.247626492829795
$ws.Range("H22").Value = 2.845720272745666
$ws.Range("I22").Value = 3.475253603157057
$ws.Range("J22").Value = 0.2516993557018878
$ws.Range("L22").Value = 0.2108881709999437
$ws.Range("N22").Value = 4.67379181795809
$ws.Range("C23").Value = 0.08117068124059479
$ws.Range("D23").Value = 0.04474780300229497
$ws.Range("E23").Value = 0.1138690459239289
$ws.Range("F23").Value = 4.700041032170844
$ws.Range("G23").Value = 4.191102778043842
$ws.Range("H23").Value = 2.822348962675107
$ws.Range("I23").Value = 3.438587690463521
$ws.Range("J23").Value = 0.2508630729911303
$ws.Range("L23").Value = 0.2107018048309683
$ws.Range("N23").Value = 4.523002190001307
$ws.Range("C24").Value = 0.08156330140021062
$ws.Range("D24").Value = 0.04337139788721345
$ws.Range("E24").Value = 0.1140422263295147
$ws.Range("F24").Value = 4.520851223648521
$ws.Range("G24").Value = 3.981697274730777
$ws.Range("H24").Value = 2.736725153440773
$ws.Range("I24").Value = 3.303401883963886
$ws.Range("J24").Value = 0.2479146873638882
$ws.Range("L24").Value = 0.2101501760575175
$ws.Range("N24").Value = 3.950713976768498
$ws.Range("C25").Value = 0.082043555527056
$ws.Range("D25").Value = 0.04186193773093905
$ws.Range("E25").Value = 0.1143813132759952
$ws.Range("F25").Value = 4.338049164424717
$ws.Range("G25").Value = 3.765568830016605
$ws.Range("H25").Value = 2.650445154803322
$ws.Range("I25").Value = 3.165296621382481
$ws.Range("J25").Value = 0.2451978751496497
$ws.Range("L25").Value = 0.2098832758479077
$ws.Range("N25").Value = 3.331249627311138
